$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3111.111
$ws.Range("I76").Value = 3123.077
$ws.Range("J76").Value = 3080
$ws.Range("K76").Value = 3123.077
$ws.Range("L76").Value = 3080
$ws.Range("M76").Value = -2808.077
$ws.Range("N76").Value = -3710

$ws.Range("H79").Value = 3111.111
$ws.Range("I79").Value = 3123.077
$ws.Range("J79").Value = 3080
$ws.Range("K79").Value = 3123.077
$ws.Range("L79").Value = 3080
$ws.Range("M79").Value = -2031.077
$ws.Range("N79").Value = -5264

$ws.Range("H113").Value = 4434.2646
$ws.Range("I113").Value = 3607.6924
$ws.Range("J113").Value = 4945.952
$ws.Range("K113").Value = 3607.6924
$ws.Range("L113").Value = 4945.952
$ws.Range("M113").Value = -353.6923999999999
$ws.Range("N113").Value = -11453.952

$ws.Range("H132").Value = 4648.7026
$ws.Range("I132").Value = 1337.0303
$ws.Range("K132").Value = 4011.0909
$ws.Range("M132").Value = -1481.0909

$ws.Range("H137").Value = 1129.5555
$ws.Range("I137").Value = 892.069
$ws.Range("J137").Value = 1560
$ws.Range("K137").Value = 2676.207
$ws.Range("L137").Value = 4680
$ws.Range("M137").Value = -126.2069999999999
$ws.Range("N137").Value = -9780

$ws.Range("H141").Value = 1669353.8
$ws.Range("I141").Value = 2001624.5
$ws.Range("K141").Value = 6004873.5
$ws.Range("M141").Value = -5999693.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2467.2727
$ws.Range("I2").Value = 10000
$ws.Range("J2").Value = 1714
$ws.Range("K2").Value = 10000
$ws.Range("L2").Value = 1714
$ws.Range("M2").Value = -9887
$ws.Range("N2").Value = -1940

$ws.Range("H40").Value = 13185.714
$ws.Range("J40").Value = 13185.714
$ws.Range("L40").Value = 13185.714
$ws.Range("N40").Value = -13537.714

$ws.Range("H45").Value = 2250
$ws.Range("I45").Value = 1500
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 1500
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -1123
$ws.Range("N45").Value = -3754

$ws.Range("H52").Value = 17999.334
$ws.Range("J52").Value = 17999.334
$ws.Range("L52").Value = 17999.334
$ws.Range("N52").Value = -18635.334

$ws.Range("H116").Value = 2467.2727
$ws.Range("I116").Value = 10000
$ws.Range("J116").Value = 1714
$ws.Range("K116").Value = 10000
$ws.Range("L116").Value = 1714
$ws.Range("M116").Value = -7706
$ws.Range("N116").Value = -6302

$ws.Range("H132").Value = 1728.75
$ws.Range("I132").Value = 1333.0358
$ws.Range("J132").Value = 4498.75
$ws.Range("K132").Value = 3999.1074
$ws.Range("L132").Value = 13496.25
$ws.Range("M132").Value = -1469.1074
$ws.Range("N132").Value = -18556.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2467.2727
$ws.Range("I3").Value = 10000
$ws.Range("J3").Value = 1714
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 1714
$ws.Range("M3").Value = -9886
$ws.Range("N3").Value = -1942

$ws.Range("H86").Value = 1802.2354
$ws.Range("I86").Value = 1646.4706
$ws.Range("J86").Value = 1958
$ws.Range("K86").Value = 1646.4706
$ws.Range("L86").Value = 1958
$ws.Range("M86").Value = -523.4706000000001
$ws.Range("N86").Value = -4204

$ws.Range("H89").Value = 1802.2354
$ws.Range("I89").Value = 1646.4706
$ws.Range("J89").Value = 1958
$ws.Range("K89").Value = 8232.353000000001
$ws.Range("L89").Value = 9790
$ws.Range("M89").Value = -2616.353000000001
$ws.Range("N89").Value = -21022

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1968.2
$ws.Range("I31").Value = 1579.5
$ws.Range("J31").Value = 2745.6
$ws.Range("K31").Value = 1579.5
$ws.Range("L31").Value = 2745.6
$ws.Range("M31").Value = -1284.5
$ws.Range("N31").Value = -3335.6

$ws.Range("H34").Value = 1968.2
$ws.Range("I34").Value = 1579.5
$ws.Range("J34").Value = 2745.6
$ws.Range("K34").Value = 1579.5
$ws.Range("L34").Value = 2745.6
$ws.Range("M34").Value = -1377.5
$ws.Range("N34").Value = -3149.6

$ws.Range("H58").Value = 6891.65
$ws.Range("I58").Value = 9125.615
$ws.Range("J58").Value = 2742.8572
$ws.Range("K58").Value = 9125.615
$ws.Range("L58").Value = 2742.8572
$ws.Range("M58").Value = -8922.615
$ws.Range("N58").Value = -3148.8572

$ws.Range("H136").Value = 6891.65
$ws.Range("I136").Value = 9125.615
$ws.Range("J136").Value = 2742.8572
$ws.Range("K136").Value = 27376.845
$ws.Range("L136").Value = 8228.571599999999
$ws.Range("M136").Value = -24826.845
$ws.Range("N136").Value = -13328.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 516.8
$ws.Range("I92").Value = 502
$ws.Range("J92").Value = 526.6667
$ws.Range("K92").Value = 1506
$ws.Range("L92").Value = 1580.0001
$ws.Range("M92").Value = -258
$ws.Range("N92").Value = -4076.0001

$ws.Range("H100").Value = 3389.4736
$ws.Range("J100").Value = 3389.4736
$ws.Range("L100").Value = 10168.4208
$ws.Range("N100").Value = -11790.4208

$ws.Range("H113").Value = 842372.4399999999
$ws.Range("I113").Value = 2020820.6
$ws.Range("J113").Value = 623.7619
$ws.Range("K113").Value = 6062461.800000001
$ws.Range("L113").Value = 1871.2857
$ws.Range("M113").Value = -6060291.800000001
$ws.Range("N113").Value = -6211.2857

$ws.Range("H115").Value = 2533.4285
$ws.Range("J115").Value = 2940
$ws.Range("L115").Value = 8820
$ws.Range("N115").Value = -11170

$ws.Range("H131").Value = 918.76
$ws.Range("J131").Value = 944.0417
$ws.Range("L131").Value = 2832.1251
$ws.Range("N131").Value = -12912.1251

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 53470
$ws.Range("I22").Value = 333793.34
$ws.Range("J22").Value = 909.375
$ws.Range("K22").Value = 333793.34
$ws.Range("L22").Value = 909.375
$ws.Range("M22").Value = -333498.34
$ws.Range("N22").Value = -1499.375

$ws.Range("H27").Value = 53470
$ws.Range("I27").Value = 333793.34
$ws.Range("J27").Value = 909.375
$ws.Range("K27").Value = 333793.34
$ws.Range("L27").Value = 909.375
$ws.Range("M27").Value = -333686.34
$ws.Range("N27").Value = -1123.375

$ws.Range("H46").Value = 1422.2222
$ws.Range("I46").Value = 1542.8572
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 1542.8572
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -1354.8572
$ws.Range("N46").Value = -1376

$ws.Range("H132").Value = 4612.3125
$ws.Range("I132").Value = 4226.8184
$ws.Range("K132").Value = 12680.4552
$ws.Range("M132").Value = -10150.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2880.5
$ws.Range("I132").Value = 3067.6206
$ws.Range("J132").Value = 2518.7334
$ws.Range("K132").Value = 9202.861800000001
$ws.Range("L132").Value = 7556.2002
$ws.Range("M132").Value = -6672.861800000001
$ws.Range("N132").Value = -12616.2002
